# Apply the cryptocurrency price/volume refresh described in the commit
# ("Updated cryptos list ... with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): values are plain text (e.g. "20.50", "116.00", "1.000")
# and must stay text. If assigned while the cell's number format is "General",
# Excel re-parses them as numbers and silently drops meaningful trailing zeros /
# punctuation ("20.50" -> 20.5, "116.00" -> 116, "1.000" -> 1). Force a text
# format first, assign, then restore the default style so the saved file keeps
# the original (unstyled) cell formatting.
$priceCellRefs = @(
    "D2",
    "D3",
    "D5",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D16",
    "D18",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $priceCellRefs) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.901.41'
$ws.Range("D3").Value = '1.810.73'
$ws.Range("D5").Value = '310.38'
$ws.Range("D7").Value = '0.4631'
$ws.Range("D8").Value = '0.3707'
$ws.Range("D9").Value = '0.07367'
$ws.Range("D10").Value = '0.8753'
$ws.Range("D11").Value = '20.50'
$ws.Range("D12").Value = '1.775.83'
$ws.Range("D13").Value = '5.362'
$ws.Range("D14").Value = '92.21'
$ws.Range("D16").Value = '0.07041'
$ws.Range("D18").Value = '0.000008707'
$ws.Range("D21").Value = '26.924.75'
$ws.Range("D22").Value = '5.319'
$ws.Range("D23").Value = '10.61'
$ws.Range("D24").Value = '2.069.54'
$ws.Range("D25").Value = '1.898'
$ws.Range("D26").Value = '151.62'
$ws.Range("D27").Value = '18.40'
$ws.Range("D28").Value = '2.159'
$ws.Range("D29").Value = '5.333'
$ws.Range("D30").Value = '116.00'
$ws.Range("D31").Value = '0.08918'
$ws.Range("D32").Value = '0.7543'
$ws.Range("D33").Value = '1.161'
$ws.Range("D34").Value = '4.459'
$ws.Range("D35").Value = '2.918'
$ws.Range("D37").Value = '1.104'
$ws.Range("D38").Value = '0.01972'
$ws.Range("D39").Value = '2.449'
$ws.Range("D40").Value = '0.05252'
$ws.Range("D41").Value = '2.933'
$ws.Range("D42").Value = '0.5331'
$ws.Range("D43").Value = '7.216'
$ws.Range("D44").Value = '0.1667'
$ws.Range("D45").Value = '8.510'
$ws.Range("D46").Value = '0.4990'
$ws.Range("D48").Value = '104.12'
$ws.Range("D49").Value = '1.000'
$ws.Range("D50").Value = '1.671'
$ws.Range("D51").Value = '0.06300'

foreach ($cellRef in $priceCellRefs) {
    $ws.Range($cellRef).Style = "Normal"
}

# --- Coin name (B), link (C) and 1h volume (E) columns: these never look like
# numbers (URLs, names, or percentages wrapped in spaces), so plain text
# assignment is safe and keeps the default cell style untouched.
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("E3").Value = '  -0.96%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("E5").Value = '  -0.87%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +4.08%  '
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("E12").Value = '  -2.85%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("E14").Value = '  -0.96%  '
$ws.Range("E15").Value = '  -3.15%  '
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("E21").Value = '  -1.71%  '
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("E28").Value = '  -5.77%  '
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  -5.33%  '
$ws.Range("E33").Value = '  -3.28%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E34").Value = '  -2.22%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("E35").Value = '  -1.72%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  -0.41%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E39").Value = '  +4.55%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E40").Value = '  -0.45%  '
$ws.Range("E41").Value = '  +2.03%  '
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("E45").Value = '  -2.32%  '
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("E47").Value = '  -3.22%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("E51").Value = '  -1.46%  '
